$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting while we
# overwrite numeric-looking strings (Excel would otherwise coerce "632.08" etc. to a number).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.225.16"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.777.49"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "632.08"
$ws.Range("E5").Value = "  +4.05%  "

$ws.Range("D6").Value = "166.37"
$ws.Range("E6").Value = "  +1.77%  "

$ws.Range("D7").Value = "3.777.69"
$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.92%  "

$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("E11").Value = "  +2.36%  "

$ws.Range("D12").Value = "6.78"
$ws.Range("E12").Value = "  -2.62%  "

$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "34.90"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "4.410.05"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "3.792.94"
$ws.Range("E16").Value = "  -1.02%  "

$ws.Range("D17").Value = "69.205.14"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").Value = "17.64"
$ws.Range("E18").Value = "  -2.29%  "

$ws.Range("E19").Value = "  +0.63%  "

$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").Value = "463.86"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "9.56"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").Value = "0.709"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").Value = "0.0000146"
$ws.Range("E24").Value = "  -1.01%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "82.90"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").Value = "12.10"
$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("D27").Value = "2.15"
$ws.Range("E27").Value = "  +2.05%  "

$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +1.42%  "

$ws.Range("D30").Value = "3.926.14"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("E32").Value = "  +3.01%  "

$ws.Range("D33").Value = "7.12"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("D34").Value = "28.53"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("E35").Value = "  +15.11%  "

$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").Value = "3.729.12"
$ws.Range("E37").Value = "  -0.40%  "

$ws.Range("D38").Value = "8.98"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").Value = "3.32"
$ws.Range("E40").Value = "  +4.33%  "

$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("E42").Value = "  -1.59%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").Value = "158.13"
$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  +5.30%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.42"
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("E48").Value = "  +0.68%  "

$ws.Range("D49").Value = "0.296"
$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").Value = "46.71"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("E51").Value = "  +0.13%  "

# Restore default (Normal) style on the price column so no stray text format sticks.
$priceRange.Style = "Normal"